# se agrego un for
# Adds a new data row (iteration 2) to the "MCP-INT-01" sheet, which also
# becomes the active/selected sheet (moving tabSelected away from
# "MCP-INT-02" automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCP-INT-01")

# Activating this sheet makes it the workbook's active tab and marks its
# sheetView as tabSelected, while clearing tabSelected on the previously
# active sheet (MCP-INT-02).
$ws.Activate()

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "MCP-INT-01"
$ws.Range("C3").Value = "Chrome"
$ws.Range("D3").Value = "https://www.coppel.com/"
$ws.Range("E3").Value = "Julian"
$ws.Range("F3").Value = "Medina"
$ws.Range("G3").Value = "Masculino"
$ws.Range("H3").Value = 6672108838
$ws.Range("I3").Value = "HolaMundo"
$ws.Range("J3").Value = "No"

# Final selection on the sheet, per the saved view state.
$ws.Range("B7").Select()
